$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Total" (sheet1.xml): append daily rows 424-429
# -----------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("Total")

$totalRows = @(
    @{ Row=424; A=44391; B=12995; E=6; F=0; G=322; H=12941 },
    @{ Row=425; A=44392; B=13011; E=6; F=0; G=322; H=12957 },
    @{ Row=426; A=44393; B=13017; E=4; F=0; G=323; H=13000 },
    @{ Row=427; A=44394; B=13017; E=4; F=0; G=323; H=13000 },
    @{ Row=428; A=44395; B=13017; E=4; F=0; G=323; H=13000 },
    @{ Row=429; A=44396; B=13035; E=3; F=0; G=323; H=13000 }
)

foreach ($r in $totalRows) {
    $row = $r.Row

    $cellA = $wsTotal.Range("A$row")
    $cellA.Value = $r.A
    $cellA.NumberFormat = "d-mmm"
    $cellA.HorizontalAlignment = -4108

    $cellB = $wsTotal.Range("B$row")
    $cellB.Value = $r.B
    $cellB.HorizontalAlignment = -4108

    $cellC = $wsTotal.Range("C$row")
    $cellC.Formula = "=B$row-H$row-G$row"
    $cellC.HorizontalAlignment = -4108

    $cellD = $wsTotal.Range("D$row")
    $cellD.Formula = "=B$row-H$row"
    $cellD.HorizontalAlignment = -4108

    $cellE = $wsTotal.Range("E$row")
    $cellE.Value = $r.E
    $cellE.HorizontalAlignment = -4108

    $cellF = $wsTotal.Range("F$row")
    $cellF.Value = $r.F
    $cellF.HorizontalAlignment = -4108

    $cellG = $wsTotal.Range("G$row")
    $cellG.Value = $r.G
    $cellG.HorizontalAlignment = -4108

    $cellH = $wsTotal.Range("H$row")
    $cellH.Value = $r.H
    $cellH.HorizontalAlignment = -4108

    $prev = $row - 1

    $cellI = $wsTotal.Range("I$row")
    $cellI.Formula = "=B$row-B$prev"
    $cellI.HorizontalAlignment = -4108

    $cellJ = $wsTotal.Range("J$row")
    $cellJ.Formula = "=H$row-H$prev"
    $cellJ.HorizontalAlignment = -4108

    $cellK = $wsTotal.Range("K$row")
    $cellK.Formula = "=G$row-G$prev"
    $cellK.HorizontalAlignment = -4108

    $cellL = $wsTotal.Range("L$row")
    $cellL.Formula = "=E$row-E$prev"
    $cellL.HorizontalAlignment = -4108

    $cellM = $wsTotal.Range("M$row")
    $cellM.Formula = "=F$row-F$prev"
    $cellM.HorizontalAlignment = -4108
}

# -----------------------------------------------------------------
# Sheet "Dados Hoje" (sheet2.xml): append daily rows 3-8
# -----------------------------------------------------------------
$wsHoje = $wb.Worksheets.Item("Dados Hoje")

$hojeRows = @(
    @{ Row=3; A=44391; B=672;  C=32563; D=196; E=12995; F=19372; G=12941; H=6; I=0; J=322 },
    @{ Row=4; A=44392; B=596;  C=32640; D=180; E=13011; F=19449; G=12957; H=6; I=0; J=322 },
    @{ Row=5; A=44393; B=566;  C=32701; D=141; E=13017; F=19543; G=13000; H=4; I=0; J=323 },
    @{ Row=6; A=44394; B=566;  C=32701; D=141; E=13017; F=19543; G=13000; H=4; I=0; J=323 },
    @{ Row=7; A=44395; B=566;  C=32701; D=141; E=13017; F=19543; G=13000; H=4; I=0; J=323 },
    @{ Row=8; A=44396; B=511;  C=32766; D=57;  E=13035; F=19674; G=13000; H=3; I=0; J=323 }
)

foreach ($r in $hojeRows) {
    $row = $r.Row

    $cellA = $wsHoje.Range("A$row")
    $cellA.Value = $r.A
    $cellA.NumberFormat = "d-mmm"
    $cellA.HorizontalAlignment = -4108

    foreach ($col in @("B","C","D","E","F","G","H","I","J")) {
        $cell = $wsHoje.Range("$col$row")
        $cell.Value = $r[$col]
        $cell.HorizontalAlignment = -4108
    }
}

# Ghost / leftover empty-formatted cells seen past the data block
$wsHoje.Range("L8").Borders.LineStyle = -4142
$wsHoje.Range("K10").Borders.LineStyle = -4142
$wsHoje.Range("G13").Borders.LineStyle = -4142

# -----------------------------------------------------------------
# Selections / active cells
# -----------------------------------------------------------------
$wsTotal.Activate()
$excel.Goto($wsTotal.Range("I429"), $false)

$wsHoje.Activate()
$excel.Goto($wsHoje.Range("J8"), $false)
